# Commit: "Implement trust-based isolation and BCA data import"
#
# The only substantive content edit captured by the target diff is a
# worksheet rename: the first sheet ("Sem-3", the active/tab-selected
# sheet) becomes "Sem2".
#
# (Everything else that differs in the canonical-XML diff --
# fileVersion/rupBuild, the xr:revisionPtr documentId GUID, the
# bookViews window geometry, and the defaultRowHeight / x14ac:dyDescent
# / best-fit column-width churn sprinkled through sheetFormatPr and
# <cols> on both sheets -- is the mechanical fallout of the workbook
# having been resaved by a newer Excel build recalculating font
# metrics/default row height, not a discrete user action, so it isn't
# reproduced here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sem2"
